$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'0"
$ws.Range("B2").Value = "patients with NAFLD"
$ws.Range("C2").Value = "effects of moderate alcohol consumption"
$ws.Range("H2").Value = "does not compare NAFLD patients with a control group"
$ws.Range("A3").Value = "'1"
$ws.Range("B3").Value = "children with liver disease"
$ws.Range("E3").Value = "not specified"
$ws.Range("H3").Value = "focuses on liver disease in children, not specifically NAFLD"
$ws.Range("A4").Value = "'2"
$ws.Range("H4").Value = "does not mention NAFLD or any relevant outcomes"
$ws.Range("A5").Value = "'3"
$ws.Range("B5").Value = "patients with NAFLD"
$ws.Range("E5").Value = "extrahepatic malignancies"
$ws.Range("F5").Value = "narrative review"
$ws.Range("H5").Value = "is a review, not a cohort study, and lacks comparison data"
$ws.Range("A6").Value = "'4"
$ws.Range("C6").Value = "not specified"
$ws.Range("E6").Value = "not specified"
$ws.Range("F6").Value = "not specified"
$ws.Range("H6").Value = "does not mention NAFLD or any relevant outcomes"
$ws.Range("A7").Value = "'5"
$ws.Range("B7").Value = "not specified"
$ws.Range("C7").Value = "not specified"
$ws.Range("E7").Value = "not specified"
$ws.Range("H7").Value = "does not mention NAFLD or any relevant outcomes"
$ws.Range("A8").Value = "'6"
$ws.Range("B8").Value = "patients with NAFLD"
$ws.Range("D8").Value = "general population"
$ws.Range("E8").Value = "risk of colorectal polyps"
$ws.Range("F8").Value = "meta-analysis of observational studies"
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = "matches all PICOS criteria with relevant outcomes and comparisons"
$ws.Range("A9").Value = "'9"
$ws.Range("C9").Value = "scutellarin effects"
$ws.Range("H9").Value = "does not mention NAFLD or relevant outcomes"
$ws.Range("A10").Value = "'7"
$ws.Range("A11").Value = "'8"
$ws.Range("B11").Value = "not applicable"
$ws.Range("C11").Value = "not applicable"
$ws.Range("D11").Value = "not applicable"
$ws.Range("E11").Value = "not applicable"
$ws.Range("F11").Value = "not applicable"
$ws.Range("H11").Value = "Not processed - Empty abstract"
$ws.Range("A12").Value = "'11"
$ws.Range("B12").Value = "not specified"
$ws.Range("C12").Value = "curcumin applications for health promotion"
$ws.Range("H12").Value = "does not address NAFLD or related outcomes"
$ws.Range("A13").Value = "'12"
$ws.Range("C13").Value = "effects of intestinal fungi on health"
$ws.Range("E13").Value = "correlation with various diseases including NAFLD"
$ws.Range("H13").Value = "does not focus on NAFLD management or cancer outcomes"
$ws.Range("A14").Value = "'13"
$ws.Range("B14").Value = "patients with LMNA variants"
$ws.Range("D14").Value = "not specified"
$ws.Range("E14").Value = "various clinical manifestations"
$ws.Range("F14").Value = "not specified"
$ws.Range("H14").Value = "does not involve NAFLD or cancer outcomes"
$ws.Range("A15").Value = "'14"
$ws.Range("C15").Value = "influence of sex and gender on health"
$ws.Range("H15").Value = "does not address NAFLD or related outcomes"
$ws.Range("A16").Value = "'15"
$ws.Range("B16").Value = "premenopausal breast cancer patients"
$ws.Range("C16").Value = "comparison of TOR and TAM treatments"
$ws.Range("E16").Value = "incidence of fatty liver and other side effects"
$ws.Range("F16").Value = "prospective randomized clinical study"
$ws.Range("H16").Value = "focuses on breast cancer, not NAFLD management or cancer incidence"
$ws.Range("A17").Value = "'16"
$ws.Range("B17").Value = "pregnant women with liver diseases"
$ws.Range("C17").Value = "management of liver diseases"
$ws.Range("D17").Value = "not specified"
$ws.Range("E17").Value = "adverse maternal and fetal outcomes"
$ws.Range("F17").Value = "not specified"
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = "discusses liver diseases in pregnancy, not specifically NAFLD or cancer outcomes"
$ws.Range("A18").Value = "'17"
$ws.Range("C18").Value = "overview of F. prausnitzii features"
$ws.Range("E18").Value = "correlation with intestinal disorders"
$ws.Range("F18").Value = "not specified"
$ws.Range("H18").Value = "does not address NAFLD or related outcomes"
$ws.Range("A19").Value = "'18"
$ws.Range("B19").Value = "not specified"
$ws.Range("C19").Value = "targeting senescent cells for CKD"
$ws.Range("D19").Value = "not specified"
$ws.Range("E19").Value = "role of cellular senescence in kidney fibrosis"
$ws.Range("F19").Value = "not specified"
$ws.Range("H19").Value = "focuses on CKD, not NAFLD or cancer incidence"
$ws.Range("A20").Value = "'10"
